$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: Carry bypass Adder
$ws.Range("C9").Value = 15276.3
$ws.Range("F9").Value = 4223.7
$ws.Range("G9").Value = 78
$ws.Range("H9").Value = 207
$ws.Range("I9").Value = 10.944514
$ws.Range("J9").Value = 4.015822
$ws.Range("K9").Value = 23.751881000000001

# Row 10: Carry Select Adder
$ws.Range("C10").Value = 15613.1
$ws.Range("F10").Value = 3886.9
$ws.Range("G10").Value = 68
$ws.Range("H10").Value = 238
$ws.Range("I10").Value = 10.567257
$ws.Range("J10").Value = 4.3081699999999996
$ws.Range("K10").Value = 25.079926

# Update selection to M10
$ws.Range("M10").Select()
